$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("F2").Value = 8.25
$ws.Range("D3").Value = 10.69
$ws.Range("C4").Value = 9.31
$ws.Range("E4").Value = 9.949999999999999
$ws.Range("F4").Value = 10.25
$ws.Range("D5").Value = 10.05
$ws.Range("F5").Value = 10.3
$ws.Range("G5").Value = 9.140000000000001
$ws.Range("B6").Value = 11.75
$ws.Range("D6").Value = 9.75
$ws.Range("E6").Value = 9.699999999999999
$ws.Range("J6").Value = 8.6
$ws.Range("E7").Value = 10.86
$ws.Range("H7").Value = 9.77
$ws.Range("G8").Value = 10.23
$ws.Range("F10").Value = 11.4
